$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at the top (row 27-28), shifting existing rows 27.. down to 29..
$ws.Rows("27:28").Insert()

# Row 27: new entry, "Zafiro rojo", Primera
$ws.Range("A27").Value = 7
$ws.Range("B27").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C27").Value = 'Ñuble'
$ws.Range("D27").Value = 44453
$ws.Range("D27").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E27").Value = 16
$ws.Range("F27").Value = 100112002
$ws.Range("G27").Value = 'Pimiento'
$ws.Range("H27").Value = 'Zafiro rojo'
$ws.Range("I27").Value = 'Primera'
$ws.Range("J27").Value = 160
$ws.Range("K27").Value = 38000
$ws.Range("L27").Value = 40000
$ws.Range("M27").Value = 39000
$ws.Range("N27").Value = '$/caja 15 kilos'
$ws.Range("O27").Value = 'Región de Arica y Parinacota'
$ws.Range("P27").Value = 2600
$ws.Range("Q27").Value = 15
$ws.Range("R27").Value = 'Hortaliza'

# Row 28: new entry, "Zafiro verde", Primera
$ws.Range("A28").Value = 7
$ws.Range("B28").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C28").Value = 'Ñuble'
$ws.Range("D28").Value = 44453
$ws.Range("D28").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = 100112002
$ws.Range("G28").Value = 'Pimiento'
$ws.Range("H28").Value = 'Zafiro verde'
$ws.Range("I28").Value = 'Primera'
$ws.Range("J28").Value = 120
$ws.Range("K28").Value = 36000
$ws.Range("L28").Value = 37000
$ws.Range("M28").Value = 36500
$ws.Range("N28").Value = '$/caja 15 kilos'
$ws.Range("O28").Value = 'Región de Arica y Parinacota'
$ws.Range("P28").Value = 2433
$ws.Range("Q28").Value = 15
$ws.Range("R28").Value = 'Hortaliza'

Write-Output "done"
